$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Komentarz" (column I) notes for weeks 1-6 (rows 5-10) ---
$ws.Range("I5").Value = "Po zapoznaniu się z dokumentacją, konieczne było wyjaśnienie paru niejasności z prowadzącymi przedmiot (1h). Przedyskutowanie problemu z grupą (2h)."
$ws.Range("I6").Value = "Założenie repozytorium zgodnie z zaleceniami (2h). Stworzenie szkieletu klienta (1h)."
$ws.Range("I7").Value = "Stworzenie całego klienta (wersja niedziałająca). Ponadto opracowanie części projektu Common wspólnego dla wielu komponentów (całość 3h)."
$ws.Range("I8").Value = "Dopracowanie clienta - wersja poprawnie nawiązująca połączenie i wysyłająca cokolwiek (3h). Nadprogramowe (2h) na poprawki związane ze zmianą sposobu komunikacji."
$ws.Range("I9").Value = "Dopracowanie clienta - wersja działająca także z serwerami innych zespołów (3h). Nadprogramowe (6h) związane z ponowną zmianą sposobu komunikacji. Kolejne (4h) więcej na stworzenie dokumentacji oraz UnitTestów."
$ws.Range("I10").Value = "Ostatnie poprawki, zwłaszcza w wyglądzie kodu. Poprawienie drobnych błędów. (Całość 3h)"

# --- Row heights: every data row (5-18) grows to fit the longer comments ---
$ws.Range("A5:A18").RowHeight = 129.95

# --- Sheet view: zoom way out and scroll down so the whole calendar is visible ---
$excel.ActiveWindow.Zoom = 25
$ws.Range("P10").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

# --- Page setup: landscape, smaller scale to fit the taller rows ---
$ws.PageSetup.Orientation = 2
$ws.PageSetup.Zoom = 57
